# Applies the "Update countries & provincias Spain" data refresh to the Pais sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country labels that shifted rows because the underlying ranking re-sorted ---
$countryLabels = @{
    "A20" = "Turquia"
    "A21" = "Italia"
    "A82" = "Libano"
    "A83" = "Paraguay"
    "A84" = "Madagascar"
    "A92" = "Grecia"
    "A93" = "Malasia"
}
foreach ($addr in $countryLabels.Keys) {
    $ws.Range($addr).Value = $countryLabels[$addr]
}

# --- Refreshed case statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes) ---
$stats = @{
    "B4" = 6015317
    "C4" = 14952
    "D4" = 3316394
    "E4" = 2514964
    "G4" = 306
    "H4" = 183959
    "B5" = 3731022
    "C5" = 9018
    "E5" = 704178
    "G5" = 240
    "H5" = 117996
    "B6" = 3377908
    "C6" = 70159
    "D6" = 2577990
    "E6" = 738272
    "G6" = 1017
    "H6" = 61646
    "B13" = 404102
    "C13" = 1737
    "D13" = 377922
    "E13" = 15108
    "G13" = 82
    "H13" = 11072
    "B16" = 330368
    "C16" = 1522
    "G16" = 12
    "H16" = 41477
    "B20" = 263998
    "C20" = 1491
    "D20" = 240792
    "E20" = 16997
    "G20" = 26
    "H20" = 6209
    "B21" = 263949
    "C21" = 1411
    "D21" = 206554
    "E21" = 21932
    "G21" = 5
    "H21" = 35463
    "D51" = 55139
    "E51" = 1406
    "B82" = 14937
    "C82" = 689
    "D82" = 4133
    "E82" = 10658
    "G82" = 7
    "H82" = 146
    "B83" = 14872
    "C83" = 0
    "D83" = 8134
    "E83" = 6491
    "G83" = 0
    "H83" = 247
    "B84" = 14592
    "C84" = 38
    "D84" = 13686
    "E84" = 722
    "G84" = 3
    "H84" = 184
    "B92" = 9531
    "C92" = 251
    "D92" = 3804
    "E92" = 5473
    "G92" = 6
    "H92" = 254
    "B93" = 9296
    "C93" = 5
    "D93" = 8994
    "E93" = 177
    "H93" = 125
    "B152" = 1429
    "C152" = 18
    "E152" = 1222
    "B175" = 431
    "C175" = 1
    "E175" = 85
}
foreach ($addr in $stats.Keys) {
    $ws.Range($addr).Value = $stats[$addr]
}

# --- Update the "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 18:28"
